$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Sent By" header in column G
$ws.Cells.Item(1, 7).Value = "Sent By"

# Existing rows 2-4 get an (empty) Sent By cell since that info is unknown for them
$ws.Cells.Item(2, 7).Style = "Normal"
$ws.Cells.Item(3, 7).Style = "Normal"
$ws.Cells.Item(4, 7).Style = "Normal"

# Add new row 5 with a new SMS record that includes who it was sent by
$ws.Cells.Item(5, 1).Value = "ADM001"
$ws.Cells.Item(5, 2).Value = "Aarav Kumar Sharma"
$ws.Cells.Item(5, 3).Value = "Grade 1-Section A"
$ws.Cells.Item(5, 4).Value = "03 February 2025, 12:00 AM"
$ws.Cells.Item(5, 5).Value = "Hi Student 1, This is a test message"
$ws.Cells.Item(5, 6).Value = "Pending"
$ws.Cells.Item(5, 7).Value = "John Smith"
